$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DemoWebShop")

$ws.Range("F2").Value = "XaXzmD"
$ws.Range("G2").Value = "CFycbe"
$ws.Range("C2").Value = "CGcZNtrtt@gmail.com"
